# Generate Report for Handoff
# Adds a new tracked file (826e2da0-b989-4d2b-9dc4-3d6012052603) as row 3
# on each of the three worksheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$newFileGuid = "826e2da0-b989-4d2b-9dc4-3d6012052603"
$newFileMd   = "$newFileGuid.md"
$dateFormat  = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(3, 1).Value = $newFileMd
$wsOverview.Cells.Item(3, 2).Value = "Ready for handoff"
$wsOverview.Cells.Item(3, 3).Value = "Ready for handoff"
$wsOverview.Cells.Item(3, 4).Value = "2016-03-24 02:41:07"
$wsOverview.Cells.Item(3, 4).NumberFormat = $dateFormat

$wsOverview.Hyperlinks.Add(
    $wsOverview.Cells.Item(3, 1),
    "https://github.com/OpenLocalizationTest/oltest/blob/2a859e6940e39e2f3ac7781383507142e8509c74/e2e/$newFileMd",
    [Type]::Missing,
    [Type]::Missing,
    $newFileMd
)

# ---------------------------------------------------------------------------
# zh-cn sheet: Source File Name | File Extension | Status | Latest Handoff
# File | Latest Handoff Datetime | ... | Latest Handback DateTime | ... |
# Handoff Reason | ...
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$zhCnXlf = "$newFileGuid.813c2931c7a30731a9d33d08049df3a3ff9ff3d1.zh-cn.xlf"

$wsZhCn.Cells.Item(3, 1).Value = $newFileMd
$wsZhCn.Cells.Item(3, 2).Value = ".md"
$wsZhCn.Cells.Item(3, 3).Value = "Ready for handoff"
$wsZhCn.Cells.Item(3, 4).Value = $zhCnXlf
$wsZhCn.Cells.Item(3, 5).Value = "2016-03-24 02:41:03"
$wsZhCn.Cells.Item(3, 5).NumberFormat = $dateFormat
$wsZhCn.Cells.Item(3, 8).Value = "0001-01-01 00:00:00"
$wsZhCn.Cells.Item(3, 8).NumberFormat = $dateFormat
$wsZhCn.Cells.Item(3, 10).Value = "Include"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Cells.Item(3, 1),
    "https://github.com/OpenLocalizationTest/oltest/blob/2a859e6940e39e2f3ac7781383507142e8509c74/e2e/$newFileMd",
    [Type]::Missing,
    [Type]::Missing,
    $newFileMd
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Cells.Item(3, 4),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3c67e6b5a130a79ea196eca45d4389ac46305a08/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhCnXlf",
    [Type]::Missing,
    [Type]::Missing,
    $zhCnXlf
)

# ---------------------------------------------------------------------------
# de-de sheet: same shape as zh-cn
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$deDeXlf = "$newFileGuid.813c2931c7a30731a9d33d08049df3a3ff9ff3d1.de-de.xlf"

$wsDeDe.Cells.Item(3, 1).Value = $newFileMd
$wsDeDe.Cells.Item(3, 2).Value = ".md"
$wsDeDe.Cells.Item(3, 3).Value = "Ready for handoff"
$wsDeDe.Cells.Item(3, 4).Value = $deDeXlf
$wsDeDe.Cells.Item(3, 5).Value = "2016-03-24 02:41:07"
$wsDeDe.Cells.Item(3, 5).NumberFormat = $dateFormat
$wsDeDe.Cells.Item(3, 8).Value = "0001-01-01 00:00:00"
$wsDeDe.Cells.Item(3, 8).NumberFormat = $dateFormat
$wsDeDe.Cells.Item(3, 10).Value = "Include"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Cells.Item(3, 1),
    "https://github.com/OpenLocalizationTest/oltest/blob/2a859e6940e39e2f3ac7781383507142e8509c74/e2e/$newFileMd",
    [Type]::Missing,
    [Type]::Missing,
    $newFileMd
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Cells.Item(3, 4),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3f3c9b1606e334037032ac41f9293f606c74c6d2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deDeXlf",
    [Type]::Missing,
    [Type]::Missing,
    $deDeXlf
)
